$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 - rename "Test" to "mmr" and update metric values
$ws.Range("A2").Value = "mmr"
$ws.Range("B2").Value = 0.7134361190331808
$ws.Range("C2").Value = 0.5110410933649582
$ws.Range("D2").Value = 0.03778289221327196
$ws.Range("E2").Value = 0.1931026812795729
$ws.Range("F2").Value = 0.00909090909090909
$ws.Range("G2").Value = 0.005726110502104429
$ws.Range("H2").Value = 0.005227462901366726
$ws.Range("I2").Value = 0.008377837609220288
$ws.Range("J2").Value = 0.00231934626038035
$ws.Range("K2").Value = 0.01796536796536797
$ws.Range("L2").Value = 0.003506940450404207
$ws.Range("M2").Value = 0.03778289221327196

# Add new row 3 - "mf" model with its metric values
# (B3, C3 and N3 are intentionally left blank - no data for those metrics)
$ws.Range("A3").Value = "mf"
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("D3").Value = 0.03011123897199847
$ws.Range("E3").Value = 0.180425360411634
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0.03011123897199847
